# Natmi following Dr Hou advice
# Recomputed LR-pair statistics for Fgf2-Sdc2: the sending/target clusters now
# include "ECs" in addition to "FAPs" and "sCs", giving a full 3x3 matrix of
# sending-cluster x target-cluster combinations (rows 2-10) with updated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf2"
$ws.Range("C2").Value = "Sdc2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.6462393333333333
$ws.Range("H2").Value = 1.938718
$ws.Range("I2").Value = 0.03461850536298827
$ws.Range("J2").Value = 0.03461850536298827
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.321445333333333
$ws.Range("N2").Value = 3.964336
$ws.Range("O2").Value = 0.01021782062667047
$ws.Range("P2").Value = 0.01021782062667047
$ws.Range("Q2").Value = 0.8539699512497777
$ws.Range("R2").Value = 7.685729561248
$ws.Range("S2").Value = 0.0003537256781624438
$ws.Range("T2").Value = 0.0003537256781624438

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf2"
$ws.Range("C3").Value = "Sdc2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.6462393333333333
$ws.Range("H3").Value = 1.938718
$ws.Range("I3").Value = 0.03461850536298827
$ws.Range("J3").Value = 0.03461850536298827
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 105.9632263333333
$ws.Range("N3").Value = 317.889679
$ws.Range("O3").Value = 0.819340166699254
$ws.Range("P3").Value = 0.8193401666992541
$ws.Range("Q3").Value = 68.47760474350244
$ws.Range("R3").Value = 616.2984426915219
$ws.Range("S3").Value = 0.02836433195498983
$ws.Range("T3").Value = 0.02836433195498983

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf2"
$ws.Range("C4").Value = "Sdc2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.6462393333333333
$ws.Range("H4").Value = 1.938718
$ws.Range("I4").Value = 0.03461850536298827
$ws.Range("J4").Value = 0.03461850536298827
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.04284166666666
$ws.Range("N4").Value = 66.128525
$ws.Range("O4").Value = 0.1704420126740755
$ws.Range("P4").Value = 0.1704420126740755
$ws.Range("Q4").Value = 14.24495130343889
$ws.Range("R4").Value = 128.20456173095
$ws.Range("S4").Value = 0.005900447729835997
$ws.Range("T4").Value = 0.005900447729835996

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf2"
$ws.Range("C5").Value = "Sdc2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 15.322826
$ws.Range("H5").Value = 45.968478
$ws.Range("I5").Value = 0.8208310864042159
$ws.Range("J5").Value = 0.8208310864042158
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.321445333333333
$ws.Range("N5").Value = 3.964336
$ws.Range("O5").Value = 0.01021782062667047
$ws.Range("P5").Value = 0.01021782062667047
$ws.Range("Q5").Value = 20.24827691117867
$ws.Range("R5").Value = 182.234492200608
$ws.Range("S5").Value = 0.008387104805673326
$ws.Range("T5").Value = 0.008387104805673326

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf2"
$ws.Range("C6").Value = "Sdc2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 15.322826
$ws.Range("H6").Value = 45.968478
$ws.Range("I6").Value = 0.8208310864042159
$ws.Range("J6").Value = 0.8208310864042158
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 105.9632263333333
$ws.Range("N6").Value = 317.889679
$ws.Range("O6").Value = 0.819340166699254
$ws.Range("P6").Value = 0.8193401666992541
$ws.Range("Q6").Value = 1623.656079504285
$ws.Range("R6").Value = 14612.90471553856
$ws.Range("S6").Value = 0.67253987916636
$ws.Range("T6").Value = 0.67253987916636

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf2"
$ws.Range("C7").Value = "Sdc2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 15.322826
$ws.Range("H7").Value = 45.968478
$ws.Range("I7").Value = 0.8208310864042159
$ws.Range("J7").Value = 0.8208310864042158
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 22.04284166666666
$ws.Range("N7").Value = 66.128525
$ws.Range("O7").Value = 0.1704420126740755
$ws.Range("P7").Value = 0.1704420126740755
$ws.Range("Q7").Value = 337.7586274038833
$ws.Range("R7").Value = 3039.82764663495
$ws.Range("S7").Value = 0.1399041024321825
$ws.Range("T7").Value = 0.1399041024321825

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fgf2"
$ws.Range("C8").Value = "Sdc2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.698388
$ws.Range("H8").Value = 8.095164
$ws.Range("I8").Value = 0.1445504082327959
$ws.Range("J8").Value = 0.1445504082327959
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.321445333333333
$ws.Range("N8").Value = 3.964336
$ws.Range("O8").Value = 0.01021782062667047
$ws.Range("P8").Value = 0.01021782062667047
$ws.Range("Q8").Value = 3.565772230122667
$ws.Range("R8").Value = 32.091950071104
$ws.Range("S8").Value = 0.001476990142834698
$ws.Range("T8").Value = 0.001476990142834699

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fgf2"
$ws.Range("C9").Value = "Sdc2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.698388
$ws.Range("H9").Value = 8.095164
$ws.Range("I9").Value = 0.1445504082327959
$ws.Range("J9").Value = 0.1445504082327959
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 105.9632263333333
$ws.Range("N9").Value = 317.889679
$ws.Range("O9").Value = 0.819340166699254
$ws.Range("P9").Value = 0.8193401666992541
$ws.Range("Q9").Value = 285.9298983791507
$ws.Range("R9").Value = 2573.369085412356
$ws.Range("S9").Value = 0.1184359555779042
$ws.Range("T9").Value = 0.1184359555779042

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fgf2"
$ws.Range("C10").Value = "Sdc2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.698388
$ws.Range("H10").Value = 8.095164
$ws.Range("I10").Value = 0.1445504082327959
$ws.Range("J10").Value = 0.1445504082327959
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 22.04284166666666
$ws.Range("N10").Value = 66.128525
$ws.Range("O10").Value = 0.1704420126740755
$ws.Range("P10").Value = 0.1704420126740755
$ws.Range("Q10").Value = 59.48013943923333
$ws.Range("R10").Value = 535.3212549531
$ws.Range("S10").Value = 0.02463746251205698
$ws.Range("T10").Value = 0.02463746251205698
